$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.776.77"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "2.619.76"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.61%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").Value = "  +4.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.394"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.20%  "

$ws.Range("E11").Value = "  +1.91%  "

$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.73%  "

$ws.Range("D14").Value = "3.090.35"

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "63.707.81"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000170"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +13.54%  "

$ws.Range("D17").Value = "2.616.93"
$ws.Range("E17").Value = "  -0.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.16%  "

$ws.Range("E21").Value = "  +2.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.71%  "

$ws.Range("E24").Value = "  -3.05%  "

$ws.Range("E25").Value = "  +0.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "547.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.162"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.39%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.23%  "

$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0903"
$ws.Range("E32").Value = "  +7.42%  "

$ws.Range("E33").Value = "  +4.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.420"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "164.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.09%  "

$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "168.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.72%  "

$ws.Range("E43").Value = "  +4.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.74%  "

$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0583"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.96%  "

$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.98%  "

$ws.Range("E47").Value = "  +0.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0251"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0972"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.48%  "

$ws.Range("D51").Value = "0.0₆0232"
$ws.Range("E51").Value = "  +18.77%  "
